$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Time ratio" column (N) ---------------------------------------
$ws.Range("N1").Value = "Time ratio"

# Row 2 (single visible row in its own shared-formula block upstream)
$ws.Range("N2").Formula = "=L2/M2"

# Rows 6:11 (visible rows sharing one formula block, same as L/M columns)
$ws.Range("N6:N11").Formula = "=L6/M6"

# Rows 15:16 (visible rows sharing another formula block)
$ws.Range("N15:N16").Formula = "=L15/M15"

# Row 18 (visible, standalone formula - row 17 is hidden and skipped)
$ws.Range("N18").Formula = "=L18/M18"

# --- Summary table used for the AVERAGE() calculation (rows 26-36) -----
$ws.Range("N26").Value = 14.25
$ws.Range("N27").Value = 3.4098360655737707
$ws.Range("N28").Value = 3.043715846994536
$ws.Range("N29").Value = 8.4868735083532219
$ws.Range("N30").Value = 0.42408376963350786
$ws.Range("N31").Value = 1.173913043478261
$ws.Range("N32").Value = 8.3673469387755102
$ws.Range("N33").Value = 1.2793650793650795
$ws.Range("N34").Value = 3.2433628318584069
$ws.Range("N35").Value = 3.5761929194458695
$ws.Range("N36").Formula = "=AVERAGE(N26:N35)"

# --- Resize the chart so its bottom edge follows the extra row ---------
$co = $ws.ChartObjects(1)
$co.Height = $co.Height + 18.112913385826772

# --- View state: zoom + selection (matches the author's final screen) --
$win = $excel.ActiveWindow
$win.Zoom = 124
$ws.Range("J18").Select()
